$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the Table5 data by the "Value" column (B) instead of "Data Point" (A).
# Values are written directly (row by row) so the existing per-row
# formatting/borders stay anchored to their row position, exactly as
# Excel's table sort visually reorders the data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 20
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = 20
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 40
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 50
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 140
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = 200

# Fill in the helper values used for the Median calculation
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 40

# Fill in helper values used for the Mean calculation
$ws.Range("J4").Value = 7
$ws.Range("J3").Formula = "=SUM(Table5[Value])"
$ws.Range("J2").Formula = "=J3/J4"

# Fill in helper value used for the Mode calculation
$ws.Range("J8").Value = 20

# Update the formulas that pull the Median / Mode values into the summary box
$ws.Range("F7").Formula = "=F4"
$ws.Range("F9").Formula = "=J8"

$excel.Calculate()
